# repull data, push all data, mean calculation
# Update column F (dSF) values for the rows where the recalculated
# final score differential changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -6
    3  = -9
    4  = -6
    7  = -5
    10 = -1
    12 = 0
    15 = -2
    19 = 2
    24 = 5
    27 = -2
    28 = -4
    40 = -2
    41 = 0
    47 = -4
    50 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
